# Updated C3DC phs000466 queries
# The "Treatment" query (row 5, TreatmentTab) wrapped its REPLACE() call in a
# redundant CONCAT(...) - drop the CONCAT() wrapper so the formula reads
# REPLACE(trt.treatment_agent, ';', ', ') directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$treatmentCell = $ws.Range("B5")
$oldText = $treatmentCell.Text

$needle = "CONCAT(REPLACE(trt.treatment_agent, ';', ', '))"
$replacement = "REPLACE(trt.treatment_agent, ';', ', ')"
$newText = $oldText -replace [regex]::Escape($needle), $replacement

$treatmentCell.Value = $newText

# The B4/B5 query cells picked up a (cosmetically identical) font-size touch
# during the same editing session.
$ws.Range("B4").Font.Size = 12
$ws.Range("B5").Font.Size = 12

# Selection ended up back on B2 when the workbook was last saved.
$ws.Range("B2").Select()
